# Update the cryptos price list (Price + Volume(1h) columns) with the
# latest scraped values. Cells whose new "Price" text would otherwise be
# auto-parsed as a number by Excel are forced back to text (NumberFormat
# "@" + Style "Normal") so they keep their original inline-string shape
# (preserving values like "0.0820" / "4.00" instead of losing trailing
# zeros to numeric coercion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.923.19'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.237.25'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.93%  '
$ws.Range('E7').Value = '  -3.05%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0820'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.09%  '
$ws.Range('E13').Value = '  -2.68%  '
$ws.Range('D14').Value = '2.578.44'
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '2.235.84'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.837'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.03%  '
$ws.Range('D18').Value = '43.829.46'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.23%  '
$ws.Range('D20').Value = '0.0₃0969'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('E24').Value = '  -7.22%  '
$ws.Range('E25').Value = '  -8.28%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.28%  '
$ws.Range('E30').Value = '  -8.60%  '
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0830'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.05%  '
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('E36').Value = '  -8.09%  '
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('E38').Value = '  -3.31%  '
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.93%  '
$ws.Range('E42').Value = '  -5.79%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').Value = '1.701.93'
$ws.Range('E44').Value = '  -3.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.44'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.96%  '
$ws.Range('E46').Value = '  -6.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '71.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '56.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.67%  '
$ws.Range('E51').Value = '  -4.15%  '
